# Lab 2 Book1.xlsx - add a new "C" frequency row (row 10), pushing the
# repeated measurement rows that followed it down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 10 (shifts old rows 11-19 down to 12-20, and
# Excel automatically re-derives the fill-down formulas for the rows
# that moved).
$ws.Rows("10:10").Insert()

# Populate the new row with the "C" note (2093 * 2 doubled again -> 4186 Hz)
$ws.Range("B10").Value = "C"
$ws.Range("C10").Value = 4186
$ws.Range("D10").Formula = "=C10*2"
$ws.Range("E10").Formula = "=65536-(22118400/D10)"

# Match the saved selection/viewport from the edited workbook.
$ws.Range("E25").Select() | Out-Null
